$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6449.6665
$ws.Range("I40").Value = 6578.2856
$ws.Range("J40").Value = 5999.5
$ws.Range("K40").Value = 6578.2856
$ws.Range("L40").Value = 5999.5
$ws.Range("M40").Value = -6403.2856
$ws.Range("N40").Value = -6349.5
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369
$ws.Range("H112").Value = 1941.9333
$ws.Range("J112").Value = 1941.9333
$ws.Range("L112").Value = 5825.7999
$ws.Range("N112").Value = -8041.7999
$ws.Range("H113").Value = 13335130
$ws.Range("I113").Value = 41667972
$ws.Range("J113").Value = 2028.0588
$ws.Range("K113").Value = 41667972
$ws.Range("L113").Value = 2028.0588
$ws.Range("M113").Value = -41664718
$ws.Range("N113").Value = -8536.058800000001
$ws.Range("H132").Value = 4556.4614
$ws.Range("I132").Value = 4966.25
$ws.Range("J132").Value = 2302.625
$ws.Range("K132").Value = 14898.75
$ws.Range("L132").Value = 6907.875
$ws.Range("M132").Value = -12368.75
$ws.Range("N132").Value = -11967.875
$ws.Range("H135").Value = 933.53845
$ws.Range("I135").Value = 789.35297
$ws.Range("J135").Value = 1914
$ws.Range("K135").Value = 7104.17673
$ws.Range("L135").Value = 17226
$ws.Range("M135").Value = -4569.17673
$ws.Range("N135").Value = -22296
$ws.Range("H137").Value = 1043840.9
$ws.Range("I137").Value = 1390712.2
$ws.Range("K137").Value = 4172136.6
$ws.Range("M137").Value = -4169586.6
$ws.Range("H138").Value = 3105.4028
$ws.Range("I138").Value = 2569.8333
$ws.Range("J138").Value = 3487.9524
$ws.Range("K138").Value = 7709.499899999999
$ws.Range("L138").Value = 10463.8572
$ws.Range("M138").Value = -2569.499899999999
$ws.Range("N138").Value = -20743.8572
$ws.Range("H141").Value = 1844.5536
$ws.Range("I141").Value = 1194.62
$ws.Range("K141").Value = 3583.86
$ws.Range("M141").Value = 1596.14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3155.327
$ws.Range("I32").Value = 3155.327
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3155.327
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2868.327
$ws.Range("H61").Value = 1990.5938
$ws.Range("I61").Value = 1826.625
$ws.Range("J61").Value = 2482.5
$ws.Range("K61").Value = 1826.625
$ws.Range("L61").Value = 2482.5
$ws.Range("M61").Value = -1614.625
$ws.Range("N61").Value = -2906.5
$ws.Range("H74").Value = 232785
$ws.Range("I74").Value = 265325.28
$ws.Range("K74").Value = 265325.28
$ws.Range("M74").Value = -264451.28
$ws.Range("H77").Value = 232785
$ws.Range("I77").Value = 265325.28
$ws.Range("K77").Value = 1326626.4
$ws.Range("M77").Value = -1322258.4
$ws.Range("H88").Value = 3833.1
$ws.Range("I88").Value = 818
$ws.Range("K88").Value = 818
$ws.Range("M88").Value = -412
$ws.Range("H91").Value = 3833.1
$ws.Range("I91").Value = 818
$ws.Range("K91").Value = 818
$ws.Range("M91").Value = 586
$ws.Range("H122").Value = 2622.8333
$ws.Range("I122").Value = 2407.5386
$ws.Range("J122").Value = 3182.6
$ws.Range("K122").Value = 7222.6158
$ws.Range("L122").Value = 9547.799999999999
$ws.Range("M122").Value = -4772.6158
$ws.Range("N122").Value = -14447.8
$ws.Range("H132").Value = 5466802
$ws.Range("I132").Value = 1892.4
$ws.Range("K132").Value = 5677.200000000001
$ws.Range("M132").Value = -3147.200000000001
$ws.Range("H136").Value = 1990.5938
$ws.Range("I136").Value = 1826.625
$ws.Range("J136").Value = 2482.5
$ws.Range("K136").Value = 5479.875
$ws.Range("L136").Value = 7447.5
$ws.Range("M136").Value = -2929.875
$ws.Range("N136").Value = -12547.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2383.5254
$ws.Range("I134").Value = 2047
$ws.Range("J134").Value = 4528.875
$ws.Range("K134").Value = 6141
$ws.Range("L134").Value = 13586.625
$ws.Range("M134").Value = -3606
$ws.Range("N134").Value = -18656.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5055.965
$ws.Range("I31").Value = 4127.091
$ws.Range("K31").Value = 4127.091
$ws.Range("M31").Value = -3832.091
$ws.Range("H34").Value = 5055.965
$ws.Range("I34").Value = 4127.091
$ws.Range("K34").Value = 4127.091
$ws.Range("M34").Value = -3925.091
$ws.Range("H58").Value = 2831.8462
$ws.Range("I58").Value = 2218.3333
$ws.Range("J58").Value = 3357.7144
$ws.Range("K58").Value = 2218.3333
$ws.Range("L58").Value = 3357.7144
$ws.Range("M58").Value = -2015.3333
$ws.Range("N58").Value = -3763.7144
$ws.Range("H86").Value = 24959
$ws.Range("I86").Value = 34965
$ws.Range("K86").Value = 34965
$ws.Range("M86").Value = -33842
$ws.Range("H89").Value = 24959
$ws.Range("I89").Value = 34965
$ws.Range("K89").Value = 174825
$ws.Range("M89").Value = -169209
$ws.Range("H107").Value = 1351.6923
$ws.Range("I107").Value = 978
$ws.Range("J107").Value = 1787.6666
$ws.Range("K107").Value = 978
$ws.Range("L107").Value = 1787.6666
$ws.Range("M107").Value = 942
$ws.Range("N107").Value = -5627.6666
$ws.Range("H132").Value = 23150486
$ws.Range("I132").Value = 16668590
$ws.Range("K132").Value = 50005770
$ws.Range("M132").Value = -50003240
$ws.Range("H134").Value = 2363
$ws.Range("I134").Value = 2252.6584
$ws.Range("K134").Value = 6757.975199999999
$ws.Range("M134").Value = -4222.975199999999
$ws.Range("H136").Value = 2831.8462
$ws.Range("I136").Value = 2218.3333
$ws.Range("J136").Value = 3357.7144
$ws.Range("K136").Value = 6654.999899999999
$ws.Range("L136").Value = 10073.1432
$ws.Range("M136").Value = -4104.999899999999
$ws.Range("N136").Value = -15173.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1060.9
$ws.Range("I113").Value = 309.625
$ws.Range("J113").Value = 1334.091
$ws.Range("K113").Value = 928.875
$ws.Range("L113").Value = 4002.273
$ws.Range("M113").Value = 1241.125
$ws.Range("N113").Value = -8342.272999999999
$ws.Range("H131").Value = 9184.690000000001
$ws.Range("I131").Value = 39844.5
$ws.Range("J131").Value = 1970.6177
$ws.Range("K131").Value = 119533.5
$ws.Range("L131").Value = 5911.8531
$ws.Range("M131").Value = -114493.5
$ws.Range("N131").Value = -15991.8531
$ws.Range("H132").Value = 1599.6666
$ws.Range("I132").Value = 1599.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14396.9994
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11866.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 41668964
$ws.Range("I80").Value = 58825564
$ws.Range("J80").Value = 2929.7144
$ws.Range("K80").Value = 58825564
$ws.Range("L80").Value = 2929.7144
$ws.Range("M80").Value = -58824566
$ws.Range("N80").Value = -4925.7144
$ws.Range("H83").Value = 41668964
$ws.Range("I83").Value = 58825564
$ws.Range("J83").Value = 2929.7144
$ws.Range("K83").Value = 294127820
$ws.Range("L83").Value = 14648.572
$ws.Range("M83").Value = -294122828
$ws.Range("N83").Value = -24632.572
$ws.Range("H107").Value = 5551.4375
$ws.Range("I107").Value = 370
$ws.Range("J107").Value = 8660.299999999999
$ws.Range("K107").Value = 370
$ws.Range("L107").Value = 8660.299999999999
$ws.Range("M107").Value = 1550
$ws.Range("N107").Value = -12500.3
$ws.Range("H113").Value = 2206.5
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 2213
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 2213
$ws.Range("M113").Value = -30
$ws.Range("N113").Value = -6553
$ws.Range("H138").Value = 107999
$ws.Range("J138").Value = 107999
$ws.Range("L138").Value = 107999
$ws.Range("N138").Value = -118279
$ws.Range("H140").Value = 72752.336
$ws.Range("J140").Value = 72752.336
$ws.Range("L140").Value = 72752.336
$ws.Range("N140").Value = -83112.336
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 232.6
$ws.Range("I16").Value = 232.6
$ws.Range("K16").Value = 232.6
$ws.Range("M16").Value = -62.59999999999999
$ws.Range("H46").Value = 1779.6
$ws.Range("I46").Value = 1899.3334
$ws.Range("J46").Value = 1600
$ws.Range("K46").Value = 1899.3334
$ws.Range("L46").Value = 1600
$ws.Range("M46").Value = -1711.3334
$ws.Range("N46").Value = -1976
$ws.Range("H82").Value = 1248.0889
$ws.Range("I82").Value = 1283.8334
$ws.Range("J82").Value = 1105.1111
$ws.Range("K82").Value = 1283.8334
$ws.Range("L82").Value = 1105.1111
$ws.Range("M82").Value = -922.8334
$ws.Range("N82").Value = -1827.1111
$ws.Range("H85").Value = 1248.0889
$ws.Range("I85").Value = 1283.8334
$ws.Range("J85").Value = 1105.1111
$ws.Range("K85").Value = 1283.8334
$ws.Range("L85").Value = 1105.1111
$ws.Range("M85").Value = -35.83339999999998
$ws.Range("N85").Value = -3601.1111
$ws.Range("H132").Value = 3504.423
$ws.Range("I132").Value = 2966.7058
$ws.Range("K132").Value = 8900.117400000001
$ws.Range("M132").Value = -6370.117400000001
$ws.Range("H136").Value = 2376.775
$ws.Range("I136").Value = 2216.6765
$ws.Range("J136").Value = 3284
$ws.Range("K136").Value = 6650.029500000001
$ws.Range("L136").Value = 9852
$ws.Range("M136").Value = -4100.029500000001
$ws.Range("N136").Value = -14952
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6411911
$ws.Range("I122").Value = 1358.1666
$ws.Range("J122").Value = 27780422
$ws.Range("K122").Value = 4074.4998
$ws.Range("L122").Value = 83341266
$ws.Range("M122").Value = -1624.4998
$ws.Range("N122").Value = -83346166
$ws.Range("H132").Value = 6412762.5
$ws.Range("I132").Value = 7938650
$ws.Range("J132").Value = 4034.3
$ws.Range("K132").Value = 23815950
$ws.Range("L132").Value = 12102.9
$ws.Range("M132").Value = -23813420
$ws.Range("N132").Value = -17162.9
Write-Output "done"
